$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sentiment labels in column A
$ws.Range("A1").Value = "Very Negative"
$ws.Range("A2").Value = "Very Positive"

# Update the counts in column B
$ws.Range("B1").Value = 2
$ws.Range("B3").Value = 97
